# Sua mensagem de commit
# Mark a batch of "Pendente" incidents as "Resolvido" on both sheets.
# On the ITI sheet the resolved rows are additionally highlighted in yellow.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# SPN sheet: flip Status (column I) from "Pendente" to "Resolvido" for the
# rows that got resolved. No formatting change on this sheet.
# ----------------------------------------------------------------------
$wsSPN = $wb.Worksheets.Item("SPN")
$spnRows = @(2, 3, 4, 7, 8, 11)
foreach ($r in $spnRows) {
    $wsSPN.Range("I$r").Value = "Resolvido"
}

# ----------------------------------------------------------------------
# ITI sheet: flip Status (column I) from "Pendente" to "Resolvido" for the
# resolved rows, and highlight those cells with a yellow fill.
# ----------------------------------------------------------------------
$wsITI = $wb.Worksheets.Item("ITI")
$itiRows = @(2, 3, 4, 5, 6, 7, 8, 9, 12, 15, 16, 17, 19, 20, 22, 24, 25, 26, 27, 28, 29)
foreach ($r in $itiRows) {
    $cell = $wsITI.Range("I$r")
    $cell.Value = "Resolvido"
    $cell.Interior.Color = 65535
}

# ----------------------------------------------------------------------
# Restore the active selection on each sheet / window position.
# ----------------------------------------------------------------------
$wsSPN.Activate() | Out-Null
$wsSPN.Range("I17").Select() | Out-Null

$wsITI.Activate() | Out-Null
$wsITI.Range("I10").Select() | Out-Null

$win = $excel.ActiveWindow
$win.Left = 28680
$win.Top = -120
